$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28

# New transaction row, matching the layout/format of the existing rows above:
#   A = Data (date, dd/mm/yyyy), B = Prelevante (text), C = Materiale (text),
#   D = UnitaMisura (text), E = Quantita (number)
$ws.Cells.Item($row, 1).VerticalAlignment = -4160  # xlTop
$ws.Cells.Item($row, 1).NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item($row, 1).Value = (Get-Date -Year 2018 -Month 5 -Day 21 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

$ws.Cells.Item($row, 2).VerticalAlignment = -4160
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "Lanzi Patrizia"

$ws.Cells.Item($row, 3).VerticalAlignment = -4160
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "Aghi 3-9"

$ws.Cells.Item($row, 4).VerticalAlignment = -4160
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "Bustine"

$ws.Cells.Item($row, 5).VerticalAlignment = -4160
$ws.Cells.Item($row, 5).Value = 4
